$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "No"
$ws.Range("C1").Value = "Yes"

$ws.Range("B2").Select()
